$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row (row 4): set the Value cell (B4) to the computed fshGenerated name
$ws.Range("B4").Value = "FonctionqualifieeVs"

# "Date" row (row 8): refresh the generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
